# Fill in the measured data for "evaluarProgDinámica" (row 22) and
# "evaluarPow" (row 24), which were previously left blank ("Completar").
# Commit message: "Se implementan pow y prog dinamica" -> these two
# algorithms' dev/test metrics are now filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Métricas")

# --- Row 22: evaluarProgDinámica ---
$ws.Range("F22").Value = 12
$ws.Range("G22").Value = 0.003472222222222222
$ws.Range("H22").Value = 0.9874999999999999
$ws.Range("I22").Value = 0.9902777777777777
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 0.0020833333333333333
$ws.Range("M22").Value = 12

# --- Row 24: evaluarPow ---
$ws.Range("F24").Value = 10
$ws.Range("G24").Value = 0.003472222222222222
$ws.Range("H24").Value = 0.9812500000000001
$ws.Range("I24").Value = 0.9826388888888888
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 7

# Row 24's J/N cells had no formula at all before (row was fully blank);
# give them the same formulas used by the sibling rows in the table.
$ws.Range("J24").Formula = '=IFERROR(IF(OR(ISBLANK(H24),ISBLANK(I24)),"",IF(I24>=H24,I24-H24,"Error")),"Error")'
$ws.Range("N24").Formula = '=IFERROR(IF(OR(J24="",ISBLANK(L24)),"",J24+L24),"Error")'

$excel.Calculate()

# Move the active selection as recorded in the saved workbook.
$ws.Range("F23").Select()
